$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order/content for rows 2-16 (row 1 header, row 8/13/17 unchanged, reordering of player rows)
$data = @(
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Christian Braun", "SG,SF", "Denver Nuggets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
